$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date in column C for rows 2-6
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45174
}
